$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" updates ---
$ws1.Range("F3").Value2 = 88
$ws1.Range("F4").Value2 = 263
$ws1.Range("F6").Value2 = 539
$ws1.Range("F8").Value2 = 1998
$ws1.Range("F10").Value2 = 96
$ws1.Range("F11").Value2 = 4270
$ws1.Range("F13").Value2 = 280
$ws1.Range("F15").Value2 = 99
$ws1.Range("F16").Value2 = 20
$ws1.Range("C17").Value2 = '宜春·逆光ZERO动漫游戏展'
$ws1.Range("D17").Value2 = '市府北路10号  红林大酒店'
$ws1.Range("F17").Value2 = 14
$ws1.Range("H17").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91866'
$ws1.Range("I17").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/i1gKVM991726717574907.jpeg'
$ws1.Range("C18").Value2 = '抚州·逆光ZERO动漫游戏展'
$ws1.Range("D18").Value2 = '王安石大道2466号 保利华章希尔顿逸林酒店'
$ws1.Range("E18").Value2 = '2024.10.03 10:00-10.03 17:00'
$ws1.Range("F18").Value2 = 62
$ws1.Range("G18").Value2 = 40
$ws1.Range("H18").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91865'
$ws1.Range("I18").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/1p0DrTb91725280390796.jpeg'
$ws1.Range("C19").Value2 = '江西·JMG（江西广电）第二届UP动漫游戏博览会'
$ws1.Range("D19").Value2 = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws1.Range("E19").Value2 = '2024.10.03 09:00-10.05 18:00'
$ws1.Range("F19").Value2 = 3019
$ws1.Range("G19").Value2 = 70
$ws1.Range("H19").Value2 = 'https://show.bilibili.com/platform/detail.html?id=90599'
$ws1.Range("I19").Value2 = '//i0.hdslb.com/bfs/openplatform/202408/oZpM885D1724642687206.png'
$ws1.Range("C20").Value2 = '萍乡·AU10秋至国漫展'
$ws1.Range("D20").Value2 = '建设西路钻石公寓西南侧60米 智博篮球馆'
$ws1.Range("E20").Value2 = '2024.10.03 10:00-10.03 17:00'
$ws1.Range("F20").Value2 = 58
$ws1.Range("G20").Value2 = 45
$ws1.Range("H20").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92178'
$ws1.Range("I20").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/YWK3l0Zx1725857595232.jpeg'
$ws1.Range("C21").Value2 = '赣州·第五届半夏动漫展'
$ws1.Range("D21").Value2 = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws1.Range("E21").Value2 = '2024.10.03 10:00-10.05 17:00'
$ws1.Range("F21").Value2 = 437
$ws1.Range("G21").Value2 = 49.5
$ws1.Range("H21").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91719'
$ws1.Range("I21").Value2 = '//i2.hdslb.com/bfs/openplatform/202409/KLxI6RZQ1725270195248.jpeg'
$ws1.Range("B22").Value2 = '2024-10-03'
$ws1.Range("C22").Value2 = '鹰潭·夜穹动漫游戏嘉年华'
$ws1.Range("D22").Value2 = '南站路66号 回禾酒店（鹰潭火车站南站路店）'
$ws1.Range("E22").Value2 = '2024.10.03 10:00-10.03 17:00'
$ws1.Range("F22").Value2 = 17
$ws1.Range("G22").Value2 = 45
$ws1.Range("H22").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91960'
$ws1.Range("I22").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/EudXOPTz1725362358018.jpeg'
$ws1.Range("C23").Value2 = '九江·无限喵国潮动漫节'
$ws1.Range("D23").Value2 = '洪垅大道 智汇欣体育中心'
$ws1.Range("E23").Value2 = '2024.10.04 10:00-10.04 17:00'
$ws1.Range("F23").Value2 = 15
$ws1.Range("G23").Value2 = 30
$ws1.Range("H23").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92370'
$ws1.Range("I23").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/Mebp8k9u1725460545541.jpeg'
$ws1.Range("C24").Value2 = '宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华'
$ws1.Range("D24").Value2 = '宜春国际商贸城会展中心 宜春国际商贸城会展中心'
$ws1.Range("E24").Value2 = '2024.10.04 10:00-10.05 17:00'
$ws1.Range("F24").Value2 = 67
$ws1.Range("G24").Value2 = 55
$ws1.Range("H24").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91115'
$ws1.Range("I24").Value2 = '//i0.hdslb.com/bfs/openplatform/202408/8TOnPvxz1723000627660.jpeg'
$ws1.Range("F25").Value2 = 69
$ws1.Range("F26").Value2 = 7
$ws1.Range("F28").Value2 = 47
$ws1.Range("F30").Value2 = 8
$ws1.Range("F31").Value2 = 424
$ws1.Range("F32").Value2 = 1680
$ws1.Range("F33").Value2 = 245

# --- Sheet "全部类型" updates ---
$ws4.Range("F3").Value2 = 88
$ws4.Range("F4").Value2 = 263
$ws4.Range("F6").Value2 = 539
$ws4.Range("F8").Value2 = 1998
$ws4.Range("F10").Value2 = 96
$ws4.Range("F11").Value2 = 4270
$ws4.Range("F13").Value2 = 280
$ws4.Range("F15").Value2 = 99
$ws4.Range("F16").Value2 = 20
$ws4.Range("C17").Value2 = '宜春·逆光ZERO动漫游戏展'
$ws4.Range("D17").Value2 = '市府北路10号  红林大酒店'
$ws4.Range("F17").Value2 = 14
$ws4.Range("H17").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91866'
$ws4.Range("I17").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/i1gKVM991726717574907.jpeg'
$ws4.Range("C18").Value2 = '抚州·逆光ZERO动漫游戏展'
$ws4.Range("D18").Value2 = '王安石大道2466号 保利华章希尔顿逸林酒店'
$ws4.Range("E18").Value2 = '2024.10.03 10:00-10.03 17:00'
$ws4.Range("F18").Value2 = 62
$ws4.Range("G18").Value2 = 40
$ws4.Range("H18").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91865'
$ws4.Range("I18").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/1p0DrTb91725280390796.jpeg'
$ws4.Range("C19").Value2 = '江西·JMG（江西广电）第二届UP动漫游戏博览会'
$ws4.Range("D19").Value2 = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws4.Range("E19").Value2 = '2024.10.03 09:00-10.05 18:00'
$ws4.Range("F19").Value2 = 3019
$ws4.Range("G19").Value2 = 70
$ws4.Range("H19").Value2 = 'https://show.bilibili.com/platform/detail.html?id=90599'
$ws4.Range("I19").Value2 = '//i0.hdslb.com/bfs/openplatform/202408/oZpM885D1724642687206.png'
$ws4.Range("C20").Value2 = '萍乡·AU10秋至国漫展'
$ws4.Range("D20").Value2 = '建设西路钻石公寓西南侧60米 智博篮球馆'
$ws4.Range("E20").Value2 = '2024.10.03 10:00-10.03 17:00'
$ws4.Range("F20").Value2 = 58
$ws4.Range("G20").Value2 = 45
$ws4.Range("H20").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92178'
$ws4.Range("I20").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/YWK3l0Zx1725857595232.jpeg'
$ws4.Range("C21").Value2 = '赣州·第五届半夏动漫展'
$ws4.Range("D21").Value2 = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws4.Range("E21").Value2 = '2024.10.03 10:00-10.05 17:00'
$ws4.Range("F21").Value2 = 437
$ws4.Range("G21").Value2 = 49.5
$ws4.Range("H21").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91719'
$ws4.Range("I21").Value2 = '//i2.hdslb.com/bfs/openplatform/202409/KLxI6RZQ1725270195248.jpeg'
$ws4.Range("B22").Value2 = '2024-10-03'
$ws4.Range("C22").Value2 = '鹰潭·夜穹动漫游戏嘉年华'
$ws4.Range("D22").Value2 = '南站路66号 回禾酒店（鹰潭火车站南站路店）'
$ws4.Range("E22").Value2 = '2024.10.03 10:00-10.03 17:00'
$ws4.Range("F22").Value2 = 17
$ws4.Range("G22").Value2 = 45
$ws4.Range("H22").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91960'
$ws4.Range("I22").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/EudXOPTz1725362358018.jpeg'
$ws4.Range("C23").Value2 = '九江·无限喵国潮动漫节'
$ws4.Range("D23").Value2 = '洪垅大道 智汇欣体育中心'
$ws4.Range("E23").Value2 = '2024.10.04 10:00-10.04 17:00'
$ws4.Range("F23").Value2 = 15
$ws4.Range("G23").Value2 = 30
$ws4.Range("H23").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92370'
$ws4.Range("I23").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/Mebp8k9u1725460545541.jpeg'
$ws4.Range("C24").Value2 = '宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华'
$ws4.Range("D24").Value2 = '宜春国际商贸城会展中心 宜春国际商贸城会展中心'
$ws4.Range("E24").Value2 = '2024.10.04 10:00-10.05 17:00'
$ws4.Range("F24").Value2 = 67
$ws4.Range("G24").Value2 = 55
$ws4.Range("H24").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91115'
$ws4.Range("I24").Value2 = '//i0.hdslb.com/bfs/openplatform/202408/8TOnPvxz1723000627660.jpeg'
$ws4.Range("F25").Value2 = 69
$ws4.Range("F26").Value2 = 7
$ws4.Range("F28").Value2 = 47
$ws4.Range("F30").Value2 = 8
$ws4.Range("F31").Value2 = 425
$ws4.Range("F32").Value2 = 1680
$ws4.Range("F33").Value2 = 245
